$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("118500","17898936","6004","Royal Profit",""),
    @("118518","17898937","6020","AFSNAPKSVGAC",""),
    @("118498","17898938","1001","Current Account",""),
    @("118452","17898939","1001","Current Account",""),
    @("118518","17898940","6005","AlfalahKifayat",""),
    @("118448","17898941","1047","AsanDigCurAc",""),
    @("118452","17898942","1150","CA AKK",""),
    @("118448","17898943","1068","FreeLancerCurAc",""),
    @("118448","17898944","1005","CD BBA",""),
    @("118500","17898948","6004","Royal Profit","10"),
    @("118518","17898949","6020","AFSNAPKSVGAC","9"),
    @("118498","17898950","1001","Current Account","8"),
    @("118500","17898951","6004","Royal Profit","10"),
    @("118518","17898952","6020","AFSNAPKSVGAC","9"),
    @("118498","17898953","1001","Current Account","8"),
    @("118500","17898954","6004","Royal Profit","10"),
    @("118518","17898955","6020","AFSNAPKSVGAC","9"),
    @("118498","17898956","1001","Current Account","8"),
    @("118500","17898967","6004","Royal Profit","10"),
    @("118498","17898968","1001","Current Account","8"),
    @("118452","17898969","1001","Current Account","7"),
    @("118518","17898970","6005","AlfalahKifayat","6"),
    @("118448","17898971","1047","AsanDigCurAc","5"),
    @("118452","17898972","1150","CA AKK","2"),
    @("118448","17898973","1068","FreeLancerCurAc","1"),
    @("118448","17898975","1005","CD BBA","4"),
    @("118500","17898992","6004","Royal Profit","10"),
    @("118498","17898994","1001","Current Account","8"),
    @("118452","17898995","1001","Current Account","7"),
    @("118518","17898996","6005","AlfalahKifayat","6"),
    @("118448","17898997","1047","AsanDigCurAc","5"),
    @("118452","17898998","1150","CA AKK","2"),
    @("118448","17898999","1068","FreeLancerCurAc","1"),
    @("118448","17899000","1005","CD BBA","4"),
    @("118500","17899005","6004","Royal Profit","10"),
    @("118498","17899006","1001","Current Account","8"),
    @("118452","17899007","1001","Current Account","7"),
    @("118518","17899008","6005","AlfalahKifayat","6"),
    @("118448","17899009","1047","AsanDigCurAc","5"),
    @("118448","17899010","1068","FreeLancerCurAc","1"),
    @("118500","17899013","6004","Royal Profit","10"),
    @("118498","17899014","1001","Current Account","8"),
    @("118452","17899015","1001","Current Account","7"),
    @("118518","17899016","6005","AlfalahKifayat","6"),
    @("118448","17899017","1047","AsanDigCurAc","5"),
    @("118448","17899018","1068","FreeLancerCurAc","1"),
    @("118448","17899019","1005","CD BBA","4"),
    @("118452","17899020","1150","CA AKK","2"),
    @("118500","17899025","6004","Royal Profit","10"),
    @("118498","17899026","1001","Current Account","8"),
    @("118452","17899027","1001","Current Account","7"),
    @("118518","17899028","6005","AlfalahKifayat","6"),
    @("118452","17899029","1150","CA AKK","2"),
    @("118452","17899031","1001","Current Account","7"),
    @("118518","17899032","6005","AlfalahKifayat","6"),
    @("118448","17899033","1005","CD BBA","4"),
    @("118452","17899034","1150","CA AKK","2"),
    @("118500","17899035","6004","Royal Profit","10"),
    @("118452","17899036","1001","Current Account","7"),
    @("118518","17899037","6005","AlfalahKifayat","6"),
    @("118448","17899038","1005","CD BBA","4"),
    @("118452","17899039","1150","CA AKK","2"),
    @("118452","17899041","1001","Current Account","7"),
    @("118518","17899043","6005","AlfalahKifayat","6"),
    @("118448","17899044","1005","CD BBA","4"),
    @("118452","17899045","1150","CA AKK","2"),
    @("118500","17899046","6004","Royal Profit","10"),
    @("118498","17899047","1001","Current Account","8"),
    @("118452","17899049","1001","Current Account","7"),
    @("118448","17899050","1005","CD BBA","4"),
    @("118452","17899051","1150","CA AKK","2"),
    @("118498","17899053","1001","Current Account","8"),
    @("118452","17899054","1001","Current Account","7"),
    @("118518","17899055","6005","AlfalahKifayat","6"),
    @("118448","17899056","1005","CD BBA","4"),
    @("118452","17899057","1150","CA AKK","2"),
    @("118500","17899058","6004","Royal Profit","10"),
    @("118498","17899059","1001","Current Account","8"),
    @("118452","17899060","1001","Current Account","7"),
    @("118518","17899061","6005","AlfalahKifayat","6"),
    @("118448","17899062","1005","CD BBA","4"),
    @("118452","17899063","1150","CA AKK","2"),
    @("118498","17899064","1001","Current Account","8"),
    @("118452","17899065","1001","Current Account","7"),
    @("118518","17899067","6005","AlfalahKifayat","6"),
    @("118448","17899068","1005","CD BBA","4"),
    @("118452","17899069","1150","CA AKK","2"),
    @("118500","17899070","6004","Royal Profit","10"),
    @("118498","17899071","1001","Current Account","8"),
    @("118452","17899072","1001","Current Account","7"),
    @("118518","17899073","6005","AlfalahKifayat","6"),
    @("118448","17899074","1005","CD BBA","4")
)

$startRow = 290
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $arr = New-Object 'object[,]' 1,3
    $arr[0,0] = "'" + $row[0]
    $arr[0,1] = "'" + $row[1]
    $arr[0,2] = "'" + $row[2]
    $ws.Range("A" + $r + ":C" + $r).Value = $arr
    $ws.Cells.Item($r, 4).Value = $row[3]
    if ($row[4] -ne "") {
        $ws.Cells.Item($r, 5).Value = "'" + $row[4]
        $ws.Cells.Item($r, 5).Style = "Normal"
    }
}

$lastRow = $startRow + $data.Length - 1
$ws.Range("A" + $startRow + ":C" + $lastRow).Style = "Normal"
